# Apply the "pub runs" workbook update: refresh workbook metadata, add three
# new run records (rows 66-68), and push the blank separator / summary /
# trailing format rows down to make room for them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Insert 3 new rows just above the old blank separator row (row 66),
#    inheriting formatting from the row above (row 65) automatically.
#    Not every column had a populated cell in row 65 (e.g. I/K/L/M), so
#    re-assert the data-row font (size 8) across the full A:P block to
#    make sure every new cell - even ones we are about to fill in -
#    picks up the same cell style used throughout the table.
# ---------------------------------------------------------------------
$ws.Rows.Item(66).Resize(3).Insert()
$ws.Range("A66:P68").Font.Size = 8

# ---------------------------------------------------------------------
# 2. Row 66 - Navigation pub run (22 June 2022)
# ---------------------------------------------------------------------
$ws.Cells.Item(66, 1).Value2 = 44734
$ws.Cells.Item(66, 2).Value = "The Navigation"
$ws.Cells.Item(66, 3).Value = "Sawley"
$ws.Cells.Item(66, 4).Value = "start/end at pub"
$ws.Cells.Item(66, 5).Value = 2.45
$ws.Cells.Item(66, 6).Value = 0.029409722222222223
$ws.Cells.Item(66, 7).Formula = "=F66/E66"
$ws.Cells.Item(66, 8).Value = 1
$ws.Cells.Item(66, 9).Value = 1
$ws.Cells.Item(66, 10).Value = 1
$ws.Cells.Item(66, 14).Value = 1
$ws.Cells.Item(66, 15).Value = "Over the golf course"
$ws.Cells.Item(66, 16).Formula = "=SUM(H66:N66)*E66"

# ---------------------------------------------------------------------
# 3. Row 67 - Sitwell Arms pub run (7 Sept 2022)
# ---------------------------------------------------------------------
$ws.Cells.Item(67, 1).Value2 = 44811
$ws.Cells.Item(67, 2).Value = "The Sitwell Arms"
$ws.Cells.Item(67, 3).Value = "Horsley Woodhouse"
$ws.Cells.Item(67, 4).Value = "start/end at pub"
$ws.Cells.Item(67, 5).Value = 2.26
$ws.Cells.Item(67, 6).Value = 0.024594907407407409
$ws.Cells.Item(67, 7).Formula = "=F67/E67"
$ws.Cells.Item(67, 8).Value = 1
$ws.Cells.Item(67, 10).Value = 1
$ws.Cells.Item(67, 14).Value = 1
$ws.Cells.Item(67, 15).Value = "Refurbed pub"
$ws.Cells.Item(67, 16).Formula = "=SUM(H67:N67)*E67"

# ---------------------------------------------------------------------
# 4. Row 68 - Belper House pub run (16 Nov 2022)
# ---------------------------------------------------------------------
$ws.Cells.Item(68, 1).Value2 = 44881
$ws.Cells.Item(68, 2).Value = "Belper House"
$ws.Cells.Item(68, 3).Value = "Belper"
$ws.Cells.Item(68, 4).Value = "start/end at pub"
$ws.Cells.Item(68, 5).Value = 2.02
$ws.Cells.Item(68, 6).Value = 0.033379629629629634
$ws.Cells.Item(68, 7).Formula = "=F68/E68"
$ws.Cells.Item(68, 8).Value = 1
$ws.Cells.Item(68, 9).Value = 1
$ws.Cells.Item(68, 10).Value = 1
$ws.Cells.Item(68, 14).Value = 1
$ws.Cells.Item(68, 15).Value = "Hilly run, lively boozer"
$ws.Cells.Item(68, 16).Formula = "=SUM(H68:N68)*E68"

# ---------------------------------------------------------------------
# 5. Blank separator row is now row 69 (was 66) - already blank/formatted
#    from the Insert() shift, nothing further required.
#
#    Summary row is now row 70 (was 67) - extend its aggregate formulas to
#    cover the new data rows (3:68 instead of 3:65 / 4:68 instead of 4:65).
# ---------------------------------------------------------------------
$ws.Cells.Item(70, 5).Formula = "=SUM(E3:E68)"
$ws.Cells.Item(70, 7).Formula = "=AVERAGE(G4:G68)"
$ws.Cells.Item(70, 8).Formula = "=SUM(H3:H68)"
$ws.Cells.Item(70, 9).Formula = "=SUM(I3:I64)"
$ws.Cells.Item(70, 10).Formula = "=SUM(J3:J64)"
$ws.Cells.Item(70, 11).Formula = "=SUM(K3:K64)"
$ws.Cells.Item(70, 12).Formula = "=SUM(L3:L64)"
$ws.Cells.Item(70, 13).Formula = "=SUM(M3:M64)"
$ws.Cells.Item(70, 14).Formula = "=SUM(N3:N68)"
$ws.Cells.Item(70, 16).Formula = "=SUM(P3:P68)"

# ---------------------------------------------------------------------
# 6. Sheet view tweaks matching the author's refreshed session.
# ---------------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 39
$ws.Range("C66").Select()

$wb.Save()
